# Append the 11/05/2025 allocation row (row 65) to Sheet1, mirroring the
# "ran on 2025-11-05" daily profit-file update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 65

# Column A holds the date as literal text (matches the existing rows, which
# are all inline/shared strings like "11/04/2025", not real Excel dates).
# A bare "11/05/2025" would be auto-recognized as a date serial, so prefix
# with an apostrophe to force text, then reset the style back to Normal so
# we don't leave a stray quote-prefix format on the cell.
$ws.Cells.Item($row, 1).Value = "'11/05/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.2116545376383344
$ws.Cells.Item($row, 3).Value = 0.7883454623616656
